$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $text) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws "D2" "39.165.52"
Set-TextCell $ws "E2" "  -2.18%  "

# Row 3
Set-TextCell $ws "D3" "2.199.08"
Set-TextCell $ws "E3" "  -5.68%  "

# Row 4
Set-TextCell $ws "E4" "  +0.08%  "

# Row 5
Set-TextCell $ws "D5" "295.08"
Set-TextCell $ws "E5" "  -4.05%  "

# Row 6
Set-TextCell $ws "D6" "80.94"
Set-TextCell $ws "E6" "  -4.87%  "

# Row 7
Set-TextCell $ws "E7" "  -3.62%  "

# Row 8
Set-TextCell $ws "E8" "  +0.05%  "

# Row 9
Set-TextCell $ws "D9" "0.467"
Set-TextCell $ws "E9" "  -3.47%  "

# Row 10
Set-TextCell $ws "E10" "  -6.17%  "

# Row 11
Set-TextCell $ws "D11" "29.02"
Set-TextCell $ws "E11" "  -3.52%  "

# Row 12
Set-TextCell $ws "E12" "  -11.18%  "

# Row 13
Set-TextCell $ws "E13" "  -2.61%  "

# Row 14
Set-TextCell $ws "D14" "2.541.40"
Set-TextCell $ws "E14" "  -5.63%  "

# Row 15
Set-TextCell $ws "E15" "  -3.36%  "

# Row 16
Set-TextCell $ws "D16" "13.91"
Set-TextCell $ws "E16" "  -5.54%  "

# Row 17
Set-TextCell $ws "D17" "2.202.30"
Set-TextCell $ws "E17" "  -5.80%  "

# Row 18
Set-TextCell $ws "D18" "0.709"
Set-TextCell $ws "E18" "  -5.81%  "

# Row 19
Set-TextCell $ws "D19" "39.078.61"
Set-TextCell $ws "E19" "  -2.33%  "

# Row 20
Set-TextCell $ws "D20" "0.0₃0867"
Set-TextCell $ws "E20" "  -4.00%  "

# Row 21
Set-TextCell $ws "E21" "  -6.43%  "

# Row 22
Set-TextCell $ws "D22" "64.64"
Set-TextCell $ws "E22" "  -4.37%  "

# Row 23
Set-TextCell $ws "D23" "10.24"
Set-TextCell $ws "E23" "  -4.05%  "

# Row 24
Set-TextCell $ws "D24" "224.87"
Set-TextCell $ws "E24" "  -4.41%  "

# Row 25
Set-TextCell $ws "E25" "  -0.07%  "

# Row 26
Set-TextCell $ws "E26" "  -6.51%  "

# Row 27
Set-TextCell $ws "E27" "  -0.85%  "

# Row 28
Set-TextCell $ws "E28" "  -3.87%  "

# Row 29
Set-TextCell $ws "E29" "  +0.87%  "

# Row 30
Set-TextCell $ws "E30" "  -2.09%  "

# Row 31
Set-TextCell $ws "D31" "148.58"
Set-TextCell $ws "E31" "  -2.41%  "

# Row 32
Set-TextCell $ws "D32" "31.57"
Set-TextCell $ws "E32" "  -9.71%  "

# Row 33
Set-TextCell $ws "D33" "0.999"
Set-TextCell $ws "E33" "  -0.16%  "

# Row 34
Set-TextCell $ws "D34" "4.78"
Set-TextCell $ws "E34" "  -6.80%  "

# Row 35
Set-TextCell $ws "B35" "Hedera"
Set-TextCell $ws "C35" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell $ws "D35" "0.0693"
Set-TextCell $ws "E35" "  -4.27%  "

# Row 36
Set-TextCell $ws "B36" "WEMIXToken"
Set-TextCell $ws "C36" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws "D36" "2.33"
Set-TextCell $ws "E36" "  -4.65%  "

# Row 37
Set-TextCell $ws "D37" "0.110"
Set-TextCell $ws "E37" "  -3.56%  "

# Row 38
Set-TextCell $ws "D38" "15.21"
Set-TextCell $ws "E38" "  -3.34%  "

# Row 39
Set-TextCell $ws "D39" "0.0953"
Set-TextCell $ws "E39" "  -4.53%  "

# Row 40
Set-TextCell $ws "D40" "2.62"
Set-TextCell $ws "E40" "  -4.76%  "

# Row 41
Set-TextCell $ws "D41" "1.64"
Set-TextCell $ws "E41" "  -3.85%  "

# Row 42
Set-TextCell $ws "E42" "  -5.64%  "

# Row 43
Set-TextCell $ws "D43" "1.897.32"
Set-TextCell $ws "E43" "  -2.44%  "

# Row 44
Set-TextCell $ws "D44" "2.08"
Set-TextCell $ws "E44" "  -9.03%  "

# Row 45
Set-TextCell $ws "D45" "0.0259"
Set-TextCell $ws "E45" "  -3.23%  "

# Row 46
Set-TextCell $ws "D46" "9.01"
Set-TextCell $ws "E46" "  -2.63%  "

# Row 47
Set-TextCell $ws "D47" "15.99"
Set-TextCell $ws "E47" "  -9.29%  "

# Row 48
Set-TextCell $ws "D48" "2.60"
Set-TextCell $ws "E48" "  -3.03%  "

# Row 49
Set-TextCell $ws "D49" "71.41"
Set-TextCell $ws "E49" "  +0.09%  "

# Row 50
Set-TextCell $ws "D50" "2.408.44"
Set-TextCell $ws "E50" "  -6.13%  "

# Row 51
Set-TextCell $ws "D51" "87.03"
Set-TextCell $ws "E51" "  -6.22%  "
